{"js": "// Replace the date and each multiplication expression's operands with the\n// updated values, per the commit diff. All target strings are unique within\n// the document, so a straightforward search + replace per pair is safe.\nconst replacements = [\n  [\"2024-06-07 Friday\", \"2024-06-08 Saturday\"],\n  [\"857\u00d78=\", \"613\u00d76=\"],\n  [\"330\u00d79=\", \"922\u00d76=\"],\n  [\"835\u00d78=\", \"350\u00d75=\"],\n  [\"881\u00d74=\", \"385\u00d79=\"],\n  [\"934\u00d77=\", \"851\u00d78=\"],\n  [\"330\u00d77=\", \"342\u00d78=\"],\n  [\"603\u00d73=\", \"531\u00d78=\"],\n  [\"987\u00d73=\", \"757\u00d77=\"],\n  [\"647\u00d78=\", \"492\u00d76=\"],\n  [\"612\u00d78=\", \"538\u00d76=\"],\n  [\"482\u00d78=\", \"118\u00d72=\"],\n  [\"151\u00d73=\", \"450\u00d72=\"],\n  [\"336\u00d72=\", \"284\u00d73=\"],\n  [\"581\u00d72=\", \"721\u00d76=\"],\n  [\"969\u00d78=\", \"390\u00d79=\"],\n  [\"946\u00d75=\", \"192\u00d78=\"],\n  [\"187\u00d79=\", \"316\u00d76=\"],\n  [\"449\u00d76=\", \"351\u00d76=\"],\n  [\"398\u00d79=\", \"998\u00d74=\"],\n  [\"691\u00d74=\", \"105\u00d79=\"],\n  [\"202\u00d76=\", \"494\u00d74=\"],\n  [\"880\u00d73=\", \"774\u00d76=\"],\n  [\"765\u00d79=\", \"814\u00d77=\"],\n  [\"344\u00d75=\", \"372\u00d75=\"],\n  [\"394\u00d75=\", \"186\u00d78=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date and each multiplication expression's operands with the\n# updated values, per the commit diff. All target strings are unique within\n# the document, so a straightforward Find/Replace per pair is safe.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-06-07 Friday\", \"2024-06-08 Saturday\"),\n    @(\"857\u00d78=\", \"613\u00d76=\"),\n    @(\"330\u00d79=\", \"922\u00d76=\"),\n    @(\"835\u00d78=\", \"350\u00d75=\"),\n    @(\"881\u00d74=\", \"385\u00d79=\"),\n    @(\"934\u00d77=\", \"851\u00d78=\"),\n    @(\"330\u00d77=\", \"342\u00d78=\"),\n    @(\"603\u00d73=\", \"531\u00d78=\"),\n    @(\"987\u00d73=\", \"757\u00d77=\"),\n    @(\"647\u00d78=\", \"492\u00d76=\"),\n    @(\"612\u00d78=\", \"538\u00d76=\"),\n    @(\"482\u00d78=\", \"118\u00d72=\"),\n    @(\"151\u00d73=\", \"450\u00d72=\"),\n    @(\"336\u00d72=\", \"284\u00d73=\"),\n    @(\"581\u00d72=\", \"721\u00d76=\"),\n    @(\"969\u00d78=\", \"390\u00d79=\"),\n    @(\"946\u00d75=\", \"192\u00d78=\"),\n    @(\"187\u00d79=\", \"316\u00d76=\"),\n    @(\"449\u00d76=\", \"351\u00d76=\"),\n    @(\"398\u00d79=\", \"998\u00d74=\"),\n    @(\"691\u00d74=\", \"105\u00d79=\"),\n    @(\"202\u00d76=\", \"494\u00d74=\"),\n    @(\"880\u00d73=\", \"774\u00d76=\"),\n    @(\"765\u00d79=\", \"814\u00d77=\"),\n    @(\"344\u00d75=\", \"372\u00d75=\"),\n    @(\"394\u00d75=\", \"186\u00d78=\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n}\n"}
